$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (source diff @@ -727)
$ws.Range("H2").Value = 205.35294
$ws.Range("I2").Value = 143.33333
$ws.Range("J2").Value = 354.2
$ws.Range("K2").Value = 143.33333
$ws.Range("L2").Value = 354.2
$ws.Range("M2").Value = -30.33332999999999
$ws.Range("N2").Value = -580.2

# Row 40 (source diff @@ -2640)
$ws.Range("H40").Value = 22223142
$ws.Range("I40").Value = 30303896
$ws.Range("J40").Value = 1074.75
$ws.Range("K40").Value = 30303896
$ws.Range("L40").Value = 1074.75
$ws.Range("M40").Value = -30303721
$ws.Range("N40").Value = -1424.75

# Row 47 (source diff @@ -2998)
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()

# Row 51 (source diff @@ -3197)
$ws.Range("H51").Value = 3275.2942
$ws.Range("I51").Value = 1968.5714
$ws.Range("J51").Value = 4190
$ws.Range("K51").Value = 1968.5714
$ws.Range("L51").Value = 4190
$ws.Range("M51").Value = -1484.5714
$ws.Range("N51").Value = -5158

# Row 86 (source diff @@ -4954)
$ws.Range("H86").Value = 3049.7632
$ws.Range("I86").Value = 1529.55
$ws.Range("J86").Value = 4738.8887
$ws.Range("K86").Value = 1529.55
$ws.Range("L86").Value = 4738.8887
$ws.Range("M86").Value = -406.55
$ws.Range("N86").Value = -6984.8887

# Row 89 (source diff @@ -5107)
$ws.Range("H89").Value = 3049.7632
$ws.Range("I89").Value = 1529.55
$ws.Range("J89").Value = 4738.8887
$ws.Range("K89").Value = 7647.75
$ws.Range("L89").Value = 23694.4435
$ws.Range("M89").Value = -2031.75
$ws.Range("N89").Value = -34926.4435

# Row 92 (source diff @@ -5260)
$ws.Range("H92").Value = 1060.6666
$ws.Range("I92").Value = 962.25
$ws.Range("J92").Value = 1257.5
$ws.Range("K92").Value = 962.25
$ws.Range("L92").Value = 1257.5
$ws.Range("M92").Value = 285.75
$ws.Range("N92").Value = -3753.5

# Row 93 (source diff @@ -5312)
$ws.Range("H93").Value = 39000
$ws.Range("J93").Value = 39000
$ws.Range("L93").Value = 39000
$ws.Range("N93").Value = -43992

# Row 113 (source diff @@ -6322)
$ws.Range("H113").Value = 1608.4445
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1622
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1622
$ws.Range("M113").Value = 1754
$ws.Range("N113").Value = -8130

# Row 116 (source diff @@ -6469)
$ws.Range("H116").Value = 3107.3572
$ws.Range("I116").Value = 3724.625
$ws.Range("J116").Value = 2284.3333
$ws.Range("K116").Value = 3724.625
$ws.Range("L116").Value = 2284.3333
$ws.Range("M116").Value = -282.625
$ws.Range("N116").Value = -9168.3333

# Row 132 (source diff @@ -7265)
$ws.Range("H132").Value = 252624.45
$ws.Range("I132").Value = 259073.64
$ws.Range("J132").Value = 1106
$ws.Range("K132").Value = 777220.92
$ws.Range("L132").Value = 3318
$ws.Range("M132").Value = -774690.92
$ws.Range("N132").Value = -8378

# Row 138 (source diff @@ -7565)
$ws.Range("H138").Value = 3352.1887
$ws.Range("I138").Value = 5131.778
$ws.Range("J138").Value = 2988.182
$ws.Range("K138").Value = 15395.334
$ws.Range("L138").Value = 8964.545999999998
$ws.Range("M138").Value = -10255.334
$ws.Range("N138").Value = -19244.546

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (source diff @@ -9322)
$ws.Range("H32").Value = 3637.8315
$ws.Range("I32").Value = 3640.3777
$ws.Range("J32").Value = 3592
$ws.Range("K32").Value = 3640.3777
$ws.Range("L32").Value = 3592
$ws.Range("M32").Value = -3353.3777
$ws.Range("N32").Value = -4166

# Row 61 (source diff @@ -10731)
$ws.Range("H61").Value = 8131584
$ws.Range("I61").Value = 9525231
$ws.Range("J61").Value = 1977.6666
$ws.Range("K61").Value = 9525231
$ws.Range("L61").Value = 1977.6666
$ws.Range("M61").Value = -9525019
$ws.Range("N61").Value = -2401.6666

# Row 132 (source diff @@ -14195)
$ws.Range("H132").Value = 39617.297
$ws.Range("I132").Value = 2593.4211
$ws.Range("K132").Value = 7780.263300000001
$ws.Range("M132").Value = -5250.263300000001

# Row 136 (source diff @@ -14394)
$ws.Range("H136").Value = 8131584
$ws.Range("I136").Value = 9525231
$ws.Range("J136").Value = 1977.6666
$ws.Range("K136").Value = 28575693
$ws.Range("L136").Value = 5932.9998
$ws.Range("M136").Value = -28573143
$ws.Range("N136").Value = -11032.9998

$ws = $wb.Worksheets.Item("BSM")
# Row 62 (source diff @@ -17737)
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65 (source diff @@ -17884)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 99 (source diff @@ -19556)
$ws.Range("H99").Value = 1488.7
$ws.Range("I99").Value = 1378.5555
$ws.Range("J99").Value = 2480
$ws.Range("K99").Value = 1378.5555
$ws.Range("L99").Value = 2480
$ws.Range("M99").Value = 119.4445000000001
$ws.Range("N99").Value = -5476

# Row 105 (source diff @@ -19850)
$ws.Range("H105").Value = 2090.4119
$ws.Range("I105").Value = 2145.5
$ws.Range("J105").Value = 1833.3334
$ws.Range("K105").Value = 2145.5
$ws.Range("L105").Value = 1833.3334
$ws.Range("M105").Value = -398.5
$ws.Range("N105").Value = -5327.3334

# Row 134 (source diff @@ -21265)
$ws.Range("H134").Value = 23943
$ws.Range("I134").Value = 29925.79
$ws.Range("K134").Value = 89777.37
$ws.Range("M134").Value = -87242.37

$ws = $wb.Worksheets.Item("CRP")
# Row 99 (source diff @@ -26495)
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -4996

# Row 126 (source diff @@ -27812)
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -10940

# Row 132 (source diff @@ -28103)
$ws.Range("H132").Value = 1624.5862
$ws.Range("I132").Value = 1466.566
$ws.Range("J132").Value = 3299.6
$ws.Range("K132").Value = 4399.698
$ws.Range("L132").Value = 9898.799999999999
$ws.Range("M132").Value = -1869.698
$ws.Range("N132").Value = -14958.8

# Row 134 (source diff @@ -28204)
$ws.Range("H134").Value = 1566.8823
$ws.Range("I134").Value = 1772.16
$ws.Range("J134").Value = 996.6667
$ws.Range("K134").Value = 5316.48
$ws.Range("L134").Value = 2990.0001
$ws.Range("M134").Value = -2781.48
$ws.Range("N134").Value = -8060.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 122 (source diff @@ -34783)
$ws.Range("H122").Value = 14494317
$ws.Range("J122").Value = 2615.5
$ws.Range("L122").Value = 23539.5
$ws.Range("N122").Value = -28439.5

# Row 131 (source diff @@ -35242)
$ws.Range("H131").Value = 2287.7163
$ws.Range("I131").Value = 5380
$ws.Range("J131").Value = 1745.2106
$ws.Range("K131").Value = 16140
$ws.Range("L131").Value = 5235.6318
$ws.Range("M131").Value = -11100
$ws.Range("N131").Value = -15315.6318

# Row 132 (source diff @@ -35294)
$ws.Range("H132").Value = 111112024
$ws.Range("I132").Value = 200000740
$ws.Range("J132").Value = 1124.75
$ws.Range("K132").Value = 1800006660
$ws.Range("L132").Value = 10122.75
$ws.Range("M132").Value = -1800004130
$ws.Range("N132").Value = -15182.75

$ws = $wb.Worksheets.Item("GSM")
# Row 62 (source diff @@ -38827)
$ws.Range("H62").Value = 24146.25
$ws.Range("J62").Value = 24146.25
$ws.Range("L62").Value = 24146.25
$ws.Range("N62").Value = -25518.25

# Row 65 (source diff @@ -38968)
$ws.Range("H65").Value = 24146.25
$ws.Range("J65").Value = 24146.25
$ws.Range("L65").Value = 72438.75
$ws.Range("N65").Value = -79302.75

# Row 80 (source diff @@ -39697)
$ws.Range("H80").Value = 864647.3
$ws.Range("I80").Value = 2528377.5
$ws.Range("J80").Value = 125211.664
$ws.Range("K80").Value = 2528377.5
$ws.Range("L80").Value = 125211.664
$ws.Range("M80").Value = -2527379.5
$ws.Range("N80").Value = -127207.664

# Row 83 (source diff @@ -39847)
$ws.Range("H83").Value = 864647.3
$ws.Range("I83").Value = 2528377.5
$ws.Range("J83").Value = 125211.664
$ws.Range("K83").Value = 12641887.5
$ws.Range("L83").Value = 626058.3200000001
$ws.Range("M83").Value = -12636895.5
$ws.Range("N83").Value = -636042.3200000001

# Row 132 (source diff @@ -42218)
$ws.Range("H132").Value = 1547.0952
$ws.Range("I132").Value = 1571.075
$ws.Range("J132").Value = 1505.3914
$ws.Range("K132").Value = 4713.225
$ws.Range("L132").Value = 4516.174199999999
$ws.Range("M132").Value = -2183.225
$ws.Range("N132").Value = -9576.174199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 74 (source diff @@ -46345)
$ws.Range("H74").Value = 11889.6
$ws.Range("I74").Value = 10098.5
$ws.Range("J74").Value = 13083.667
$ws.Range("K74").Value = 10098.5
$ws.Range("L74").Value = 13083.667
$ws.Range("M74").Value = -9100.5
$ws.Range("N74").Value = -15079.667

# Row 77 (source diff @@ -46486)
$ws.Range("H77").Value = 11889.6
$ws.Range("I77").Value = 10098.5
$ws.Range("J77").Value = 13083.667
$ws.Range("K77").Value = 30295.5
$ws.Range("L77").Value = 39251.001
$ws.Range("M77").Value = -25303.5
$ws.Range("N77").Value = -49235.001

# Row 132 (source diff @@ -49145)
$ws.Range("H132").Value = 32228.805
$ws.Range("I132").Value = 36775.418
$ws.Range("J132").Value = 4039.8
$ws.Range("K132").Value = 110326.254
$ws.Range("L132").Value = 12119.4
$ws.Range("M132").Value = -107796.254
$ws.Range("N132").Value = -17179.4

# Row 136 (source diff @@ -49344)
$ws.Range("H136").Value = 8397
$ws.Range("I136").Value = 9254.166999999999
$ws.Range("J136").Value = 4968.3335
$ws.Range("K136").Value = 27762.501
$ws.Range("L136").Value = 14905.0005
$ws.Range("M136").Value = -25212.501
$ws.Range("N136").Value = -20005.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (source diff @@ -56138)
$ws.Range("H132").Value = 4886.184
$ws.Range("I132").Value = 5665.8
$ws.Range("J132").Value = 1962.625
$ws.Range("K132").Value = 16997.4
$ws.Range("L132").Value = 5887.875
$ws.Range("M132").Value = -14467.4
$ws.Range("N132").Value = -10947.875
